# Auto-generated script to apply cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.579.92'
$ws.Range('E2').Value = '  +1.10%  '
$ws.Range('D3').Value = '3.390.26'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.475'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.74'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.60%  '
$ws.Range('E10').Value = '  -0.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.386'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.29%  '
$ws.Range('D12').Value = '3.969.94'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.46'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.125'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').Value = '3.380.52'
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '61.573.08'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('E18').Value = '  -0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '391.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.553'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.06%  '
$ws.Range('E24').Value = '  -0.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000113'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.195'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.54%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.27'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.68%  '
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.40'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.54%  '
$ws.Range('B32').Value = 'USDe'
$ws.Range('C32').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.30'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.93'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '168.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.06'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.57%  '
$ws.Range('D37').Value = '3.425.22'
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('E38').Value = '  -1.35%  '
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.35%  '
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.32%  '
$ws.Range('D45').Value = '2.457.17'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.98'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.67'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.33%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('E50').Value = '  -3.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.206'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.67%  '
